$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 5674
$ws.Range("F3").Value = 7683
$ws.Range("F9").Value = 4472
$ws.Range("F10").Value = 1800
$ws.Range("F13").Value = 3009
$ws.Range("F16").Value = 227
$ws.Range("F17").Value = 564
$ws.Range("F18").Value = 485
$ws.Range("F19").Value = 488
$ws.Range("F20").Value = 347
$ws.Range("F22").Value = 1737
$ws.Range("F23").Value = 1268
$ws.Range("F24").Value = 107
$ws.Range("F25").Value = 1477
$ws.Range("F27").Value = 596
$ws.Range("F28").Value = 37
$ws.Range("F29").Value = 521
$ws.Range("F31").Value = 28
$ws.Range("F32").Value = 73
$ws.Range("F33").Value = 113
$ws.Range("F35").Value = 3236
$ws.Range("F36").Value = 725
$ws.Range("F37").Value = 50
$ws.Range("F38").Value = 188
$ws.Range("F40").Value = 1223

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 28

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5674
$ws.Range("F3").Value = 7683
$ws.Range("F9").Value = 4472
$ws.Range("F10").Value = 1800
$ws.Range("F13").Value = 3009
$ws.Range("F16").Value = 227
$ws.Range("F17").Value = 564
$ws.Range("F18").Value = 485
$ws.Range("F19").Value = 488
$ws.Range("F21").Value = 347
$ws.Range("F23").Value = 1737
$ws.Range("F24").Value = 1268
$ws.Range("F25").Value = 107
$ws.Range("F26").Value = 1477
$ws.Range("F28").Value = 596
$ws.Range("F29").Value = 37
$ws.Range("F30").Value = 521
$ws.Range("F32").Value = 28
$ws.Range("F33").Value = 73
$ws.Range("F34").Value = 113
$ws.Range("F36").Value = 3236
$ws.Range("F37").Value = 28
$ws.Range("F38").Value = 725
$ws.Range("F39").Value = 50
$ws.Range("F40").Value = 188
$ws.Range("F42").Value = 1223
